# Adding the changes we made on may 9th
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 2 (its data is superseded); all rows below shift up by one.
$ws.Rows.Item(2).Delete()

# Append the new rows of accelerometer data at the bottom (rows 21-31 after the shift).
$ws.Cells.Item(21, 1).Value = 9.343793869018556
$ws.Cells.Item(21, 2).Value = 9.33461856842041
$ws.Cells.Item(21, 3).Value = 8.671773910522461

$ws.Cells.Item(22, 1).Value = 16.0465145111084
$ws.Cells.Item(22, 2).Value = -45.62339782714844
$ws.Cells.Item(22, 3).Value = 13.98112106323242

$ws.Cells.Item(23, 1).Value = 18.59218597412109
$ws.Cells.Item(23, 2).Value = 6.586655139923096
$ws.Cells.Item(23, 3).Value = -14.55167388916016

$ws.Cells.Item(24, 1).Value = -13.51971435546875
$ws.Cells.Item(24, 2).Value = -18.57223892211914
$ws.Cells.Item(24, 3).Value = -24.33260536193848

$ws.Cells.Item(25, 1).Value = -12.54771614074707
$ws.Cells.Item(25, 2).Value = -18.61154365539551
$ws.Cells.Item(25, 3).Value = 28.26399230957031

$ws.Cells.Item(26, 1).Value = -1.716351509094239
$ws.Cells.Item(26, 2).Value = 5.821096897125244
$ws.Cells.Item(26, 3).Value = -9.188434600830078

$ws.Cells.Item(27, 1).Value = 41.91740417480469
$ws.Cells.Item(27, 2).Value = -71.66004180908203
$ws.Cells.Item(27, 3).Value = 11.51219272613525

$ws.Cells.Item(28, 1).Value = 22.04729652404785
$ws.Cells.Item(28, 2).Value = 0.830233097076416
$ws.Cells.Item(28, 3).Value = -27.47162818908692

$ws.Cells.Item(29, 1).Value = -15.56076526641846
$ws.Cells.Item(29, 2).Value = -8.609291076660156
$ws.Cells.Item(29, 3).Value = -21.82845687866211

$ws.Cells.Item(30, 1).Value = -10.87422180175781
$ws.Cells.Item(30, 2).Value = -26.33984756469727
$ws.Cells.Item(30, 3).Value = 9.715606689453123

$ws.Cells.Item(31, 1).Value = 4.591959953308105
$ws.Cells.Item(31, 2).Value = 5.822259902954102
$ws.Cells.Item(31, 3).Value = -12.29100227355957
